# Update cryptos list - GitHub Actions style refresh of prices / volume(1h) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the value to be stored as text even when it looks numeric
    # (e.g. "213.00", "18.50") so trailing zeros / exact formatting survive,
    # then restore the default "Normal" style so no stray formatting is left
    # on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.307.83"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.608.47"
$ws.Range("E3").Value = "  +0.31%  "

# Row 5 - BNB
Set-TextValue "D5" "213.00"
$ws.Range("E5").Value = "  +0.02%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.06%  "

# Row 7 - XRP
Set-TextValue "D7" "0.487"
$ws.Range("E7").Value = "  +0.14%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.69%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.10%  "

# Row 10 - Solana
Set-TextValue "D10" "18.50"
$ws.Range("E10").Value = "  +2.50%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0814"
$ws.Range("E11").Value = "  -0.24%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.832.38"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.597.71"
$ws.Range("E13").Value = "  -0.44%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.48%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.85%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "26.282.64"
$ws.Range("E16").Value = "  +0.55%  "

# Row 17 - Litecoin
$ws.Range("E17").Value = "  +2.69%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +0.74%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.00%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "201.36"
$ws.Range("E20").Value = "  -1.24%  "

# Row 21 - Uniswap
Set-TextValue "D21" "4.26"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +0.47%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +0.46%  "

# Row 24 - Toncoin
Set-TextValue "D24" "1.88"
$ws.Range("E24").Value = "  +0.84%  "

# Row 25 - Monero
Set-TextValue "D25" "143.35"
$ws.Range("E25").Value = "  +1.12%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.03%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -1.42%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.41%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +2.31%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.0497"
$ws.Range("E30").Value = "  +5.37%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.13%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.20"
$ws.Range("E32").Value = "  +2.95%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -1.45%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.19%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +1.41%  "

# Row 36 - Maker
Set-TextValue "D36" "1.163.04"
$ws.Range("E36").Value = "  +3.59%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.0168"
$ws.Range("E37").Value = "  +1.47%  "

# Row 38 - PaxDollar
$ws.Range("E38").Value = "  -0.02%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  +0.91%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +0.48%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  +1.02%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  +4.22%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "0.785"
$ws.Range("E43").Value = "  +0.23%  "

# Row 44 - RocketPoolETH
Set-TextValue "D44" "1.742.94"
$ws.Range("E44").Value = "  +0.13%  "

# Row 45 - Quant
Set-TextValue "D45" "92.02"
$ws.Range("E45").Value = "  -0.78%  "

# Rows 46 and 47 swap places (RenderToken <-> BabyDogeCoin) with updated values
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D46" "0.0₆0106"
$ws.Range("E46").Value = "  +13.79%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D47" "1.54"
$ws.Range("E47").Value = "  +1.42%  "

# Row 48 - Aave
Set-TextValue "D48" "54.10"
$ws.Range("E48").Value = "  +1.04%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.32%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.14%  "
